$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @{
    2 = @{ D = 44495; L = 'Primera'; M = 100; N = 26000; O = 27000; P = 26500; Q = '$/bandeja 10 kilos'; R = 'Provincia de Limarí'; S = 2650; T = 10 }
    3 = @{ D = 44491; L = 'Primera'; M = 150; N = 25000; O = 26000; P = 25467; Q = '$/bandeja 10 kilos'; R = 'Provincia de Limarí'; S = 2547; T = 10 }
    4 = @{ D = 44517; L = 'Primera'; M = 100; N = 25000; O = 27000; P = 26000; Q = '$/bandeja 10 kilos'; R = 'Provincia de Limarí'; S = 2600; T = 10 }
    5 = @{ D = 44467; L = 'Primera'; M = 100; N = 2700; O = 2800; P = 2750; Q = '$/kilo (en caja de 15 kilos)'; R = 'Provincia de Limarí'; S = 2750; T = 1 }
    6 = @{ D = 44467; L = 'Segunda'; M = 50; N = 2500; O = 2500; P = 2500; Q = '$/kilo (en caja de 15 kilos)'; R = 'Provincia de Limarí'; S = 2500; T = 1 }
    7 = @{ D = 44488; L = 'Primera'; M = 50; N = 25000; O = 26000; P = 25600; Q = '$/bandeja 10 kilos'; R = 'Provincia de Limarí'; S = 2560; T = 10 }
    8 = @{ D = 44446; L = 'Primera'; M = 100; N = 3200; O = 3300; P = 3250; Q = '$/kilo (en caja de 15 kilos)'; R = 'Provincia del Elquí'; S = 3250; T = 1 }
    9 = @{ D = 44483; L = 'Primera'; M = 50; N = 2600; O = 2600; P = 2600; Q = '$/kilo (en caja de 15 kilos)'; R = 'Provincia de Limarí'; S = 2600; T = 1 }
    10 = @{ D = 44483; L = 'Segunda'; M = 50; N = 2400; O = 2400; P = 2400; Q = '$/kilo (en caja de 15 kilos)'; R = 'Provincia de Limarí'; S = 2400; T = 1 }
    11 = @{ D = 44469; L = 'Primera'; M = 100; N = 28000; O = 29000; P = 28500; Q = '$/bandeja 10 kilos'; R = 'Provincia de Limarí'; S = 2850; T = 10 }
    12 = @{ D = 44505; L = 'Primera'; M = 100; N = 2200; O = 2200; P = 2200; Q = '$/kilo (en caja de 15 kilos)'; R = 'Provincia de Limarí'; S = 2200; T = 1 }
    13 = @{ D = 44505; L = 'Segunda'; M = 100; N = 1800; O = 1800; P = 1800; Q = '$/kilo (en caja de 15 kilos)'; R = 'Provincia de Limarí'; S = 1800; T = 1 }
    14 = @{ D = 44530; L = 'Primera'; M = 100; N = 2000; O = 2100; P = 2050; Q = '$/kilo (en caja de 15 kilos)'; R = 'Provincia de Limarí'; S = 2050; T = 1 }
    15 = @{ D = 44461; L = 'Primera'; M = 100; N = 29000; O = 30000; P = 29500; Q = '$/bandeja 10 kilos'; R = 'Provincia de Limarí'; S = 2950; T = 10 }
    16 = @{ D = 44462; L = 'Primera'; M = 100; N = 2900; O = 3000; P = 2950; Q = '$/kilo (en caja de 15 kilos)'; R = 'Provincia de Limarí'; S = 2950; T = 1 }
    17 = @{ D = 44462; L = 'Segunda'; M = 50; N = 2600; O = 2600; P = 2600; Q = '$/kilo (en caja de 15 kilos)'; R = 'Provincia de Limarí'; S = 2600; T = 1 }
    18 = @{ D = 44484; L = 'Primera'; M = 100; N = 25000; O = 26000; P = 25500; Q = '$/bandeja 10 kilos'; R = 'Provincia de Limarí'; S = 2550; T = 10 }
    19 = @{ D = 44516; L = 'Primera'; M = 100; N = 1900; O = 2000; P = 1950; Q = '$/kilo (en caja de 15 kilos)'; R = 'Provincia de Limarí'; S = 1950; T = 1 }
    20 = @{ D = 44516; L = 'Segunda'; M = 50; N = 1700; O = 1700; P = 1700; Q = '$/kilo (en caja de 15 kilos)'; R = 'Provincia de Limarí'; S = 1700; T = 1 }
    21 = @{ D = 44511; L = 'Primera'; M = 80; N = 25000; O = 26000; P = 25375; Q = '$/bandeja 10 kilos'; R = 'Provincia de Limarí'; S = 2538; T = 10 }
    22 = @{ D = 44160; L = 'Primera'; M = 100; N = 17000; O = 18000; P = 17500; Q = '$/bandeja 8 kilos'; R = 'Provincia de Limarí'; S = 2188; T = 8 }
    23 = @{ D = 44160; L = 'Segunda'; M = 50; N = 15000; O = 15000; P = 15000; Q = '$/bandeja 8 kilos'; R = 'Provincia de Limarí'; S = 1875; T = 8 }
    24 = @{ D = 44454; L = 'Primera'; M = 100; N = 30000; O = 31000; P = 30500; Q = '$/bandeja 10 kilos'; R = 'Provincia de Limarí'; S = 3050; T = 10 }
    25 = @{ D = 44475; L = 'Primera'; M = 100; N = 29000; O = 30000; P = 29500; Q = '$/bandeja 10 kilos'; R = 'Provincia de Limarí'; S = 2950; T = 10 }
}

foreach ($r in $rowsData.Keys) {
    $row = $rowsData[$r]
    $ws.Cells.Item($r, 4).Value2  = $row.D   # D: Fecha
    $ws.Cells.Item($r, 12).Value  = $row.L   # L: Calidad
    $ws.Cells.Item($r, 13).Value  = $row.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value  = $row.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value  = $row.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value  = $row.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value  = $row.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value  = $row.R   # R: Origen
    $ws.Cells.Item($r, 19).Value  = $row.S   # S: Precio $/Kg
    $ws.Cells.Item($r, 20).Value  = $row.T   # T: Kg / unidad
}

Write-Output "Done updating rows 2-25"
